# Última versión del proyecto
# Replace the Spanish labels ("labels_ESP") in column C with the new
# Cambodian/Khmer labels ("labels_CAM").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = "ទាប, ធម្មតា, ខ្ពស់។"
$ws.Range("C3").Value2 = "ទាប, ធម្មតា, ខ្ពស់។"
$ws.Range("C1").Value2 = "labels_CAM"
$ws.Range("C4").Value2 = "ល្អណាស់, ល្អ, មធ្យម, មិនមានសុខភាពល្អ, គ្រោះថ្នាក់, ស្លាប់"
$ws.Range("C5").Value2 = "ល្អណាស់, ល្អ, មធ្យម, មិនមានសុខភាពល្អ, គ្រោះថ្នាក់, ស្លាប់"
$ws.Range("C6").Value2 = "ល្អណាស់, ល្អ, មធ្យម, មិនមានសុខភាពល្អ, គ្រោះថ្នាក់, ស្លាប់"
$ws.Range("C7").Value2 = "ល្អណាស់, ល្អ, មធ្យម, មិនមានសុខភាពល្អ, គ្រោះថ្នាក់, ស្លាប់"
$ws.Range("C8").Value2 = "ល្អណាស់, ល្អ, មធ្យម, មិនមានសុខភាពល្អ, គ្រោះថ្នាក់, ស្លាប់"
$ws.Range("C9").Value2 = "មិនបានផ្ទុក, ផ្ទុក"

# The Battery row's label cell now wraps its (longer) Khmer text.
$ws.Range("C9").WrapText = $true

# Reflect the author's final selection in the saved sheet view.
[void]$ws.Range("C13").Select()
